$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D, E, G columns retain text formatting (values are textual, e.g. "299.62", "2.24%", "11")
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Apply updated values from the refreshed symbol list
$ws.Range("D2").Value = "299.62"
$ws.Range("E2").Value = "2.24%"
$ws.Range("G2").Value = "11"
$ws.Range("D3").Value = "42.22"
$ws.Range("E3").Value = "4.67%"
$ws.Range("G3").Value = "11"
$ws.Range("D4").Value = "5.014"
$ws.Range("E4").Value = "0.25%"
$ws.Range("G4").Value = "11"
$ws.Range("D5").Value = "0.07558"
$ws.Range("E5").Value = "3.09%"
$ws.Range("G5").Value = "11"
$ws.Range("D6").Value = "1.603"
$ws.Range("E6").Value = "2.77%"
$ws.Range("G6").Value = "11"
$ws.Range("D7").Value = "0.9379"
$ws.Range("E7").Value = "1.57%"
$ws.Range("G7").Value = "11"
$ws.Range("G8").Value = "11"
$ws.Range("E9").Value = "0.89%"
$ws.Range("G9").Value = "11"
$ws.Range("D10").Value = "0.1842"
$ws.Range("E10").Value = "1.83%"
$ws.Range("G10").Value = "11"
$ws.Range("D11").Value = "0.09081"
$ws.Range("E11").Value = "3.00%"
$ws.Range("G11").Value = "11"
$ws.Range("D12").Value = "0.04173"
$ws.Range("E12").Value = "-4.91%"
$ws.Range("G12").Value = "11"
$ws.Range("D13").Value = "0.1048"
$ws.Range("E13").Value = "-0.61%"
$ws.Range("G13").Value = "11"
$ws.Range("D14").Value = "0.001283"
$ws.Range("E14").Value = "0.89%"
$ws.Range("G14").Value = "11"
$ws.Range("D15").Value = "0.005903"
$ws.Range("E15").Value = "0.75%"
$ws.Range("G15").Value = "11"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.343"
$ws.Range("E16").Value = "-0.05%"
$ws.Range("G16").Value = "11"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "4.378"
$ws.Range("E17").Value = "2.18%"
$ws.Range("G17").Value = "11"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "0.3335"
$ws.Range("E18").Value = "0.89%"
$ws.Range("G18").Value = "11"
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D19").Value = "8.401"
$ws.Range("E19").Value = "6.68%"
$ws.Range("G19").Value = "11"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "0.1409"
$ws.Range("E20").Value = "1.35%"
$ws.Range("G20").Value = "11"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "0.3297"
$ws.Range("E21").Value = "17.76%"
$ws.Range("G21").Value = "11"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "0.04110"
$ws.Range("E22").Value = "4.77%"
$ws.Range("G22").Value = "11"
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "0.001264"
$ws.Range("E23").Value = "0.11%"
$ws.Range("G23").Value = "11"
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "0.003906"
$ws.Range("E24").Value = "6.22%"
$ws.Range("G24").Value = "11"
$ws.Range("D25").Value = "0.0001268"
$ws.Range("E25").Value = "5.66%"
$ws.Range("G25").Value = "11"
$ws.Range("G26").Value = "11"
$ws.Range("G27").Value = "11"
$ws.Range("G28").Value = "11"
$ws.Range("G29").Value = "11"
$ws.Range("G30").Value = "11"
$ws.Range("G31").Value = "11"
$ws.Range("G32").Value = "11"
$ws.Range("G33").Value = "11"
$ws.Range("G34").Value = "11"
$ws.Range("G35").Value = "11"
$ws.Range("G36").Value = "11"
$ws.Range("G37").Value = "11"
$ws.Range("D38").Value = "0.02417"
$ws.Range("E38").Value = "3.26%"
$ws.Range("G38").Value = "11"
$ws.Range("D39").Value = "0.05223"
$ws.Range("E39").Value = "2.39%"
$ws.Range("G39").Value = "11"
$ws.Range("D40").Value = "0.006793"
$ws.Range("E40").Value = "14.51%"
$ws.Range("G40").Value = "11"
$ws.Range("D41").Value = "0.007700"
$ws.Range("E41").Value = "-1.96%"
$ws.Range("G41").Value = "11"
$ws.Range("D42").Value = "0.1331"
$ws.Range("E42").Value = "3.12%"
$ws.Range("G42").Value = "11"
$ws.Range("D43").Value = "0.007381"
$ws.Range("E43").Value = "-0.13%"
$ws.Range("G43").Value = "11"
$ws.Range("D44").Value = "0.007795"
$ws.Range("E44").Value = "-3.00%"
$ws.Range("G44").Value = "11"
$ws.Range("E45").Value = "2.94%"
$ws.Range("G45").Value = "11"
$ws.Range("D46").Value = "0.00006246"
$ws.Range("E46").Value = "0.13%"
$ws.Range("G46").Value = "11"
$ws.Range("E47").Value = "-0.11%"
$ws.Range("G47").Value = "11"
$ws.Range("E48").Value = "-5.14%"
$ws.Range("G48").Value = "11"
$ws.Range("E49").Value = "-0.01%"
$ws.Range("G49").Value = "11"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").Value = "-0.11%"
$ws.Range("G50").Value = "11"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").Value = "-0.11%"
$ws.Range("G51").Value = "11"
